$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial (45202 = 2023-10-03) for every
# data row (C2:C498). Bump it by one day (45203 = 2023-10-04) for all rows.
$range = $ws.Range("C2:C498")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
